# Apply updated odds figures to the "Jogos da Semana" FlashScore sheet.
# Source: commit "Atualizando o arquivo XLSX" - refreshed odds for rows 2, 3 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Jaguares de Cordoba x Chico) ---
$ws.Range("M2").Value = 1.1
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63

# --- Row 3 (Liverpool M. x Wanderers) ---
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("AJ3").Value = 41
$ws.Range("BB3").Value = 301

# --- Row 5 (Briton Ferry x Haverfordwest) ---
$ws.Range("G5").Value = 7.6
$ws.Range("H5").Value = 4.2
$ws.Range("I5").Value = 1.37
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 2.32
$ws.Range("L5").Value = 1.87
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 3.65
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2.02
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 2.9
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.72
$ws.Range("W5").Value = 20
$ws.Range("X5").Value = 55
$ws.Range("Y5").Value = 24
$ws.Range("Z5").Value = 200
$ws.Range("AA5").Value = 100
$ws.Range("AB5").Value = 80
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 8.75
$ws.Range("AE5").Value = 20
$ws.Range("AF5").Value = 100
$ws.Range("AG5").Value = 6.6
$ws.Range("AH5").Value = 6.3
$ws.Range("AI5").Value = 8.25
$ws.Range("AJ5").Value = 8.75
$ws.Range("AK5").Value = 11.5
$ws.Range("AL5").Value = 28
$ws.Range("AM5").Value = 800
$ws.Range("AN5").Value = 8.75
$ws.Range("AO5").Value = 45
$ws.Range("AP5").Value = 45
$ws.Range("AQ5").Value = 350
$ws.Range("AR5").Value = 350
$ws.Range("AT5").Value = 2.9
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 80
$ws.Range("AW5").Value = 3.15
$ws.Range("AX5").Value = 6.3
$ws.Range("AZ5").Value = 18
$ws.Range("BA5").Value = 50
